$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.218.54'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '1.660.65'
$ws.Range('E3').Value = '  -1.16%  '
$ws.Range('E4').Value = '  +0.44%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.68'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5220'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2669'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06342'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.05'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07711'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').Value = '1.677.03'
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.428'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.60%  '
$ws.Range('D14').Value = '1.888.69'
$ws.Range('E14').Value = '  -0.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5473'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.79%  '
$ws.Range('D16').Value = '0.0₅8211'
$ws.Range('E16').Value = '  -1.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.00'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').Value = '26.237.19'
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.653'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '195.25'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.14'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.085'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.008'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.18'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1243'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.232'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.21'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.94%  '
$ws.Range('E29').Value = '  -0.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05959'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.281'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.631'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.307'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.633'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9793'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.423'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.778'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5906'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01595'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.992'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8567'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.005'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('D43').Value = '1.027.86'
$ws.Range('E43').Value = '  -4.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.87'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').Value = '1.803.06'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('D46').Value = '0.0₈109'
$ws.Range('E46').Value = '  -2.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.33'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.057'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05189'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.466'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.03%  '
